$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = $ws.Range("F4").Value2

$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = $ws.Range("F3").Value2

$ws.Range("F3").Clear()

$ws.Range("D10").Select()
